# Auto-generated: apply meteocat daily-summary update (2026-02-24 19:50 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-24 19:48:24"
$ws.Range("E3").Value = "2026-02-24 19:48:26"
$ws.Range("L3").Value = "33.1 km/h - 128º 19:16 TU"
$ws.Range("O3").Value = "4.4 °C"
$ws.Range("E4").Value = "2026-02-24 19:48:28"
$ws.Range("H4").Value = "'70%"
$ws.Range("J4").Value = "1019.8 hPa"
$ws.Range("O4").Value = "13.3 °C"
$ws.Range("E5").Value = "2026-02-24 19:48:31"
$ws.Range("H5").Value = "'28%"
$ws.Range("O5").Value = "6.0 °C"
$ws.Range("E6").Value = "2026-02-24 19:48:33"
$ws.Range("J6").Value = "1019.7 hPa"
$ws.Range("E7").Value = "2026-02-24 19:48:35"
$ws.Range("E8").Value = "2026-02-24 19:48:38"
$ws.Range("E9").Value = "2026-02-24 19:48:41"
$ws.Range("H9").Value = "'80%"
$ws.Range("O9").Value = "12.0 °C"
$ws.Range("E10").Value = "2026-02-24 19:48:43"
$ws.Range("K10").Value = "14.3 MJ/m2"
$ws.Range("O10").Value = "11.5 °C"
$ws.Range("E11").Value = "2026-02-24 19:48:45"
$ws.Range("E12").Value = "2026-02-24 19:48:48"
$ws.Range("E13").Value = "2026-02-24 19:48:50"
$ws.Range("J13").Value = "1023.3 hPa"
$ws.Range("O13").Value = "6.9 °C"
$ws.Range("E14").Value = "2026-02-24 19:48:53"
$ws.Range("E15").Value = "2026-02-24 19:48:55"
$ws.Range("H15").Value = "'75%"
$ws.Range("O15").Value = "12.2 °C"
$ws.Range("E16").Value = "2026-02-24 19:48:57"
$ws.Range("H16").Value = "'18%"
$ws.Range("K16").Value = "13.2 MJ/m2"
$ws.Range("L16").Value = "22.0 km/h - 229º 19:21 TU"
$ws.Range("O16").Value = "4.2 °C"
$ws.Range("E17").Value = "2026-02-24 19:49:00"
$ws.Range("E18").Value = "2026-02-24 19:49:02"
$ws.Range("H18").Value = "'75%"
$ws.Range("E19").Value = "2026-02-24 19:49:05"
$ws.Range("H19").Value = "'51%"
$ws.Range("E20").Value = "2026-02-24 19:49:07"
$ws.Range("E21").Value = "2026-02-24 19:49:09"
$ws.Range("J21").Value = "1022.0 hPa"
$ws.Range("E22").Value = "2026-02-24 19:49:12"
$ws.Range("O22").Value = "3.6 °C"
$ws.Range("E23").Value = "2026-02-24 19:49:14"
$ws.Range("O23").Value = "4.8 °C"
$ws.Range("E24").Value = "2026-02-24 19:49:17"
$ws.Range("J24").Value = "1021.3 hPa"
$ws.Range("O24").Value = "9.9 °C"
$ws.Range("E25").Value = "2026-02-24 19:49:19"
$ws.Range("N25").Value = "3.4 °C 19:18 TU"
$ws.Range("O25").Value = "6.9 °C"
$ws.Range("E26").Value = "2026-02-24 19:49:22"
$ws.Range("H26").Value = "'42%"
$ws.Range("O26").Value = "11.9 °C"
$ws.Range("E27").Value = "2026-02-24 19:49:24"
$ws.Range("E28").Value = "2026-02-24 19:49:27"
$ws.Range("J28").Value = "1020.1 hPa"
$ws.Range("E29").Value = "2026-02-24 19:49:29"
$ws.Range("O29").Value = "10.2 °C"
$ws.Range("E30").Value = "2026-02-24 19:49:31"
$ws.Range("H30").Value = "'74%"
$ws.Range("J30").Value = "1019.8 hPa"
$ws.Range("O30").Value = "13.3 °C"
$ws.Range("E31").Value = "2026-02-24 19:49:34"
$ws.Range("J31").Value = "1019.2 hPa"
$ws.Range("O31").Value = "15.8 °C"
$ws.Range("E32").Value = "2026-02-24 19:49:36"
$ws.Range("O32").Value = "7.6 °C"
$ws.Range("E33").Value = "2026-02-24 19:49:39"
$ws.Range("E34").Value = "2026-02-24 19:49:41"
$ws.Range("O34").Value = "4.9 °C"
$ws.Range("E35").Value = "2026-02-24 19:49:43"
$ws.Range("J35").Value = "1020.5 hPa"
$ws.Range("E36").Value = "2026-02-24 19:49:46"
$ws.Range("E37").Value = "2026-02-24 19:49:48"
$ws.Range("E38").Value = "2026-02-24 19:49:51"
$ws.Range("H38").Value = "'71%"
$ws.Range("O38").Value = "12.2 °C"
$ws.Range("E39").Value = "2026-02-24 19:49:53"
$ws.Range("H39").Value = "'35%"
$ws.Range("O39").Value = "4.6 °C"
$ws.Range("E40").Value = "2026-02-24 19:49:55"
$ws.Range("J40").Value = "1022.7 hPa"
$ws.Range("E41").Value = "2026-02-24 19:49:58"
$ws.Range("H41").Value = "'78%"
$ws.Range("J41").Value = "1020.6 hPa"
$ws.Range("E42").Value = "2026-02-24 19:50:00"
$ws.Range("H42").Value = "'86%"
$ws.Range("E43").Value = "2026-02-24 19:50:02"
$ws.Range("E44").Value = "2026-02-24 19:50:05"
$ws.Range("H44").Value = "'39%"
$ws.Range("K44").Value = "15.4 MJ/m2"
$ws.Range("E45").Value = "2026-02-24 19:50:07"
$ws.Range("O45").Value = "10.2 °C"
$ws.Range("E46").Value = "2026-02-24 19:50:10"
$ws.Range("J46").Value = "1021.3 hPa"
$ws.Range("O46").Value = "10.6 °C"
